# update scripts wuth new tpm
# Re-run of the NATMI ligand-receptor scoring for Clcf1-Crlf1 using the
# updated TPM-based expression values. This refreshes the ligand/receptor
# expression stats and the derived edge-weight / specificity columns
# (E:J and M:T) for every sending/target cluster combination on the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.1594223333333333
$ws.Range("H2").Value2 = 0.478267
$ws.Range("I2").Value2 = 0.01552338951653915
$ws.Range("J2").Value2 = 0.01552338951653915
$ws.Range("M2").Value2 = 0.1694406666666667
$ws.Range("N2").Value2 = 0.5083219999999999
$ws.Range("O2").Value2 = 0.004610192448566767
$ws.Range("P2").Value2 = 0.004610192448566767
$ws.Range("Q2").Value2 = 0.02701262644155555
$ws.Range("R2").Value2 = 0.243113637974
$ws.Range("S2").Value2 = 0.0000715658131253093
$ws.Range("T2").Value2 = 0.00007156581312530931
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 0.6666666666666666
$ws.Range("G3").Value2 = 0.1594223333333333
$ws.Range("H3").Value2 = 0.478267
$ws.Range("I3").Value2 = 0.01552338951653915
$ws.Range("J3").Value2 = 0.01552338951653915
$ws.Range("N3").Value2 = 48.51573500000001
$ws.Range("O3").Value2 = 0.4400102201629409
$ws.Range("P3").Value2 = 0.4400102201629409
$ws.Range("Q3").Value2 = 2.578163892360556
$ws.Range("R3").Value2 = 23.203475031245
$ws.Range("S3").Value2 = 0.00683045003884748
$ws.Range("T3").Value2 = 0.00683045003884748
$ws.Range("E4").Value2 = 2
$ws.Range("F4").Value2 = 0.6666666666666666
$ws.Range("G4").Value2 = 0.1594223333333333
$ws.Range("H4").Value2 = 0.478267
$ws.Range("I4").Value2 = 0.01552338951653915
$ws.Range("J4").Value2 = 0.01552338951653915
$ws.Range("M4").Value2 = 20.41213866666667
$ws.Range("N4").Value2 = 61.23641600000001
$ws.Range("O4").Value2 = 0.5553795873884922
$ws.Range("P4").Value2 = 0.5553795873884922
$ws.Range("Q4").Value2 = 3.254150774563556
$ws.Range("R4").Value2 = 29.287356971072
$ws.Range("S4").Value2 = 0.008621373664566358
$ws.Range("T4").Value2 = 0.00862137366456636
$ws.Range("I5").Value2 = 0.1862883666449807
$ws.Range("J5").Value2 = 0.1862883666449807
$ws.Range("M5").Value2 = 0.1694406666666667
$ws.Range("N5").Value2 = 0.5083219999999999
$ws.Range("O5").Value2 = 0.004610192448566767
$ws.Range("P5").Value2 = 0.004610192448566767
$ws.Range("Q5").Value2 = 0.3241649031113333
$ws.Range("R5").Value2 = 2.917484128002
$ws.Range("S5").Value2 = 0.0008588252211625272
$ws.Range("T5").Value2 = 0.0008588252211625272
$ws.Range("I6").Value2 = 0.1862883666449807
$ws.Range("J6").Value2 = 0.1862883666449807
$ws.Range("N6").Value2 = 48.51573500000001
$ws.Range("O6").Value2 = 0.4400102201629409
$ws.Range("P6").Value2 = 0.4400102201629409
$ws.Range("Q6").Value2 = 30.93924428934834
$ws.Range("S6").Value2 = 0.08196878522125262
$ws.Range("T6").Value2 = 0.0819687852212526
$ws.Range("I7").Value2 = 0.1862883666449807
$ws.Range("J7").Value2 = 0.1862883666449807
$ws.Range("M7").Value2 = 20.41213866666667
$ws.Range("N7").Value2 = 61.23641600000001
$ws.Range("O7").Value2 = 0.5553795873884922
$ws.Range("P7").Value2 = 0.5553795873884922
$ws.Range("Q7").Value2 = 39.05142185371734
$ws.Range("R7").Value2 = 351.462796683456
$ws.Range("S7").Value2 = 0.1034607562025655
$ws.Range("T7").Value2 = 0.1034607562025655
$ws.Range("G8").Value2 = 8.197245333333333
$ws.Range("H8").Value2 = 24.591736
$ws.Range("I8").Value2 = 0.7981882438384801
$ws.Range("J8").Value2 = 0.7981882438384801
$ws.Range("M8").Value2 = 0.1694406666666667
$ws.Range("N8").Value2 = 0.5083219999999999
$ws.Range("O8").Value2 = 0.004610192448566767
$ws.Range("P8").Value2 = 0.004610192448566767
$ws.Range("Q8").Value2 = 1.388946714110222
$ws.Range("R8").Value2 = 12.500520426992
$ws.Range("S8").Value2 = 0.00367980141427893
$ws.Range("T8").Value2 = 0.00367980141427893
$ws.Range("G9").Value2 = 8.197245333333333
$ws.Range("H9").Value2 = 24.591736
$ws.Range("I9").Value2 = 0.7981882438384801
$ws.Range("J9").Value2 = 0.7981882438384801
$ws.Range("N9").Value2 = 48.51573500000001
$ws.Range("O9").Value2 = 0.4400102201629409
$ws.Range("P9").Value2 = 0.4400102201629409
$ws.Range("Q9").Value2 = 132.5651274406623
$ws.Range("R9").Value2 = 1193.08614696596
$ws.Range("S9").Value2 = 0.3512109849028408
$ws.Range("T9").Value2 = 0.3512109849028408
$ws.Range("G10").Value2 = 8.197245333333333
$ws.Range("H10").Value2 = 24.591736
$ws.Range("I10").Value2 = 0.7981882438384801
$ws.Range("J10").Value2 = 0.7981882438384801
$ws.Range("M10").Value2 = 20.41213866666667
$ws.Range("N10").Value2 = 61.23641600000001
$ws.Range("O10").Value2 = 0.5553795873884922
$ws.Range("P10").Value2 = 0.5553795873884922
$ws.Range("Q10").Value2 = 167.3233084286862
$ws.Range("R10").Value2 = 1505.909775858176
$ws.Range("S10").Value2 = 0.4432974575213603
$ws.Range("T10").Value2 = 0.4432974575213603
